# binge_config.xlsx — "again change if div 0 its binge in config"
#
# The IFERROR() fallback in columns H and I changes from 0 to 2 whenever the
# corresponding divisor (F or G) is 0. H2/I2 hold their own (non-shared)
# formulas; H3:H25 / I3:I25 are a shared-formula group anchored at H3/I3, so
# rewriting the top of each range keeps the t="shared" grouping in the
# exported OOXML (only the master cell's formula text actually changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H: IFERROR(1/F,0) -> IFERROR(1/F,2)
$ws.Range("H2").Formula = "=IFERROR(1/F2,2)"
$ws.Range("H3:H25").Formula = "=IFERROR(1/F3,2)"

# Column I: IFERROR(1/G,0) -> IFERROR(1/G,2)
$ws.Range("I2").Formula = "=IFERROR(1/G2,2)"
$ws.Range("I3:I25").Formula = "=IFERROR(1/G3,2)"

# Selection moved to H2:I25 (active cell H2) and the view scrolled down so
# row 12 is at the top.
$ws.Range("H2:I25").Select()
$excel.ActiveWindow.ScrollRow = 12
